$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 2
$ws.Range("AA3").Value = 1
$ws.Range("AA5").Value = 1
$ws.Range("AA6").Value = 2
$ws.Range("AA7").Value = 2
$ws.Range("AA8").Value = 2
$ws.Range("AA9").Value = 1
$ws.Range("AA10").Value = 0
$ws.Range("AA11").Value = 1
$ws.Range("AA12").Value = 2
$ws.Range("AA13").Value = 2
$ws.Range("AA14").Value = 2
$ws.Range("AA15").Value = 2
$ws.Range("AA16").Value = 0
$ws.Range("AA19").Value = 1
$ws.Range("AA20").Value = 2
$ws.Range("AA21").Value = 2
$ws.Range("AA22").Value = 2
$ws.Range("AA23").Value = 2
$ws.Range("AA24").Value = 1
$ws.Range("AA25").Value = 1
$ws.Range("AA26").Value = 1
$ws.Range("AA27").Value = 1
$ws.Range("AA28").Value = 2
$ws.Range("AA30").Value = 1
$ws.Range("AA31").Value = 1
$ws.Range("AA32").Value = 1
$ws.Range("AA34").Value = 1
$ws.Range("AA35").Value = 0
$ws.Range("AA36").Value = 1
$ws.Range("AA37").Value = 1
$ws.Range("AA40").Value = 2
$ws.Range("AA41").Value = 2
$ws.Range("AA42").Value = 2
$ws.Range("AA43").Value = 2
$ws.Range("AA44").Value = 2
$ws.Range("AA45").Value = 1
$ws.Range("AA46").Value = 2
$ws.Range("AA47").Value = 1
$ws.Range("AA48").Value = 2
$ws.Range("AA49").Value = 2
$ws.Range("AA50").Value = 2
$ws.Range("AA52").Value = 2
$ws.Range("AA53").Value = 2
$ws.Range("AA54").Value = 1
$ws.Range("AA55").Value = 1
$ws.Range("AA56").Value = 1
$ws.Range("AA58").Value = 1
$ws.Range("AA59").Value = 1
$ws.Range("AA60").Value = 2
$ws.Range("AA61").Value = 1
$ws.Range("AA62").Value = 1
$ws.Range("AA63").Value = 2
$ws.Range("AA64").Value = 2
$ws.Range("AA65").Value = 1
$ws.Range("AA67").Value = 2
$ws.Range("AA68").Value = 1
$ws.Range("AA69").Value = 2
$ws.Range("AA70").Value = 0
$ws.Range("AA72").Value = 2
$ws.Range("AA73").Value = 2
$ws.Range("AA74").Value = 1
$ws.Range("AA75").Value = 1
$ws.Range("AA77").Value = 0
$ws.Range("AA78").Value = 1
$ws.Range("AA79").Value = 2
$ws.Range("AA80").Value = 0
$ws.Range("AA81").Value = 1
$ws.Range("AA82").Value = 2
$ws.Range("AA83").Value = 2
$ws.Range("AA85").Value = 0
$ws.Range("AA86").Value = 1
$ws.Range("AA87").Value = 1
$ws.Range("AA88").Value = 1
$ws.Range("AA89").Value = 2
$ws.Range("AA90").Value = 2
$ws.Range("AA94").Value = 2
$ws.Range("AA96").Value = 2
$ws.Range("AA98").Value = 2
$ws.Range("AA99").Value = 2
$ws.Range("AA100").Value = 1
$ws.Range("AA101").Value = 1
$ws.Range("AA102").Value = 2
$ws.Range("AA103").Value = 1
$ws.Range("AA106").Value = 1
$ws.Range("AA107").Value = 1
$ws.Range("AA108").Value = 1
$ws.Range("AA110").Value = 2
$ws.Range("AA112").Value = 2
$ws.Range("AA113").Value = 2
$ws.Range("AA115").Value = 1
$ws.Range("AA116").Value = 2
$ws.Range("AA117").Value = 1
$ws.Range("AA118").Value = 1
$ws.Range("AA119").Value = 2
$ws.Range("AA120").Value = 1
$ws.Range("AA121").Value = 1
$ws.Range("AA122").Value = 1
$ws.Range("AA123").Value = 1
$ws.Range("AA125").Value = 1
$ws.Range("AA126").Value = 2
$ws.Range("AA127").Value = 1
$ws.Range("AA128").Value = 2
$ws.Range("AA129").Value = 2
$ws.Range("AA130").Value = 2
$ws.Range("AA131").Value = 1
$ws.Range("AA133").Value = 2
$ws.Range("AA134").Value = 1
$ws.Range("AA135").Value = 1
$ws.Range("AA136").Value = 1
$ws.Range("AA139").Value = 2
$ws.Range("AA140").Value = 0
$ws.Range("AA141").Value = 0
$ws.Range("AA142").Value = 1
$ws.Range("AA144").Value = 1
$ws.Range("AA145").Value = 2
$ws.Range("AA146").Value = 2
$ws.Range("AA147").Value = 2
$ws.Range("AA148").Value = 1
$ws.Range("AA149").Value = 1
$ws.Range("AA150").Value = 1
$ws.Range("AA151").Value = 2
$ws.Range("AA152").Value = 2
$ws.Range("AA153").Value = 0
$ws.Range("AA154").Value = 1
$ws.Range("AA155").Value = 0
$ws.Range("AA157").Value = 1
$ws.Range("AA160").Value = 1
